$wb = $excel.ActiveWorkbook

# Helper values
$statusText = "Handed back: in sync with en-US"
$blueColor = 15570276  # RGB(100,149,237) == 0x6495ED, stored as BGR long for COM Font.Color

function Style-HyperlinkCell($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Underline = 2
    $range.Font.Color = $blueColor
    $range.Font.Size = 11
}

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column update (row 2 and 3)
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Target File (F) / Latest Handback File (G) -- row 2
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/63a3a75edaea2966033c16cfdf9dfb4e44d5b656/e2e/abbab7e0-a329-483e-9d20-c9deb323a86a.md", "", "", "abbab7e0-a329-483e-9d20-c9deb323a86a.md")
Style-HyperlinkCell $wsZh.Range("F2")

$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8fac1b3fe1ab662a9a982f7f6ad44e43e68f93a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/abbab7e0-a329-483e-9d20-c9deb323a86a.b74043c3bfde5436070767ea63ff90a7ff48c429.zh-cn.xlf", "", "", "abbab7e0-a329-483e-9d20-c9deb323a86a.b74043c3bfde5436070767ea63ff90a7ff48c429.zh-cn.xlf")
Style-HyperlinkCell $wsZh.Range("G2")

# Latest Target File (F) / Latest Handback File (G) -- row 3
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/63a3a75edaea2966033c16cfdf9dfb4e44d5b656/e2e/e892ccd0-7739-4a26-af12-043fadc2823d.md", "", "", "e892ccd0-7739-4a26-af12-043fadc2823d.md")
Style-HyperlinkCell $wsZh.Range("F3")

$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8fac1b3fe1ab662a9a982f7f6ad44e43e68f93a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/e892ccd0-7739-4a26-af12-043fadc2823d.e6b61b37bfdf2426bb828e2e35b716e09e313476.zh-cn.xlf", "", "", "e892ccd0-7739-4a26-af12-043fadc2823d.e6b61b37bfdf2426bb828e2e35b716e09e313476.zh-cn.xlf")
Style-HyperlinkCell $wsZh.Range("G3")

# Latest Handback DateTime (H) -- now populated with the handback timestamp
$wsZh.Range("H2").Value = "2016-03-18 16:14:47"
$wsZh.Range("H3").Value = "2016-03-18 16:14:47"

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")

# Status column update (row 2 and 3)
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Latest Target File (F) / Latest Handback File (G) -- row 2
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/63a3a75edaea2966033c16cfdf9dfb4e44d5b656/e2e/abbab7e0-a329-483e-9d20-c9deb323a86a.md", "", "", "abbab7e0-a329-483e-9d20-c9deb323a86a.md")
Style-HyperlinkCell $wsDe.Range("F2")

$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f3a5a91bd5e697f20df12db35f3a289f47b9771/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/abbab7e0-a329-483e-9d20-c9deb323a86a.b74043c3bfde5436070767ea63ff90a7ff48c429.de-de.xlf", "", "", "abbab7e0-a329-483e-9d20-c9deb323a86a.b74043c3bfde5436070767ea63ff90a7ff48c429.de-de.xlf")
Style-HyperlinkCell $wsDe.Range("G2")

# Latest Target File (F) / Latest Handback File (G) -- row 3
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/63a3a75edaea2966033c16cfdf9dfb4e44d5b656/e2e/e892ccd0-7739-4a26-af12-043fadc2823d.md", "", "", "e892ccd0-7739-4a26-af12-043fadc2823d.md")
Style-HyperlinkCell $wsDe.Range("F3")

$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f3a5a91bd5e697f20df12db35f3a289f47b9771/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/e892ccd0-7739-4a26-af12-043fadc2823d.e6b61b37bfdf2426bb828e2e35b716e09e313476.de-de.xlf", "", "", "e892ccd0-7739-4a26-af12-043fadc2823d.e6b61b37bfdf2426bb828e2e35b716e09e313476.de-de.xlf")
Style-HyperlinkCell $wsDe.Range("G3")

# Latest Handback DateTime (H) -- now populated with the handback timestamp
$wsDe.Range("H2").Value = "2016-03-18 16:14:51"
$wsDe.Range("H3").Value = "2016-03-18 16:14:51"
